# Weekly fruit/vegetable price update: insert a new latest-week row at
# row 26 (pushing the existing history rows 26-33 down to 27-34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 26:33 down by one row, inserting a blank row 26.
$ws.Rows("26").Insert()

# Populate the new row 26 with the latest weekly record.
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44504
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112032
$ws.Range("G26").Value = "Zapallo italiano"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 700
$ws.Range("K26").Value = 6500
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6750
$ws.Range("N26").Value = "$/caja 60 unidades"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 112
$ws.Range("Q26").Value = 60
$ws.Range("R26").Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Range("D26").NumberFormat = $ws.Range("D27").NumberFormat
